# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (column G) values for rows 2-10 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 3
    8  = 2
    9  = 1
    10 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
